# Update cryptocurrency price/volume data (refresh snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.411.64'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '1.841.55'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.38'
$ws.Range("E5").Value = '  -7.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5130'
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3207'
$ws.Range("E8").Value = '  -8.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06752'
$ws.Range("E9").Value = '  -1.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.70'
$ws.Range("E10").Value = '  -6.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7684'
$ws.Range("E11").Value = '  -5.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07680'
$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("D13").Value = '1.871.85'
$ws.Range("E13").Value = '  +0.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.62'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.017'
$ws.Range("E15").Value = '  -1.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.05'
$ws.Range("E17").Value = '  -1.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007886'
$ws.Range("E19").Value = '  -2.60%  '

$ws.Range("D20").Value = '26.469.14'
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").Value = '2.089.31'
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.568'
$ws.Range("E22").Value = '  -4.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.535'
$ws.Range("E23").Value = '  -5.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.949'
$ws.Range("E24").Value = '  -4.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.323'
$ws.Range("E25").Value = '  -1.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.79'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.668'
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.97'
$ws.Range("E28").Value = '  -1.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.01'
$ws.Range("E29").Value = '  +0.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.170'
$ws.Range("E30").Value = '  -4.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.153'
$ws.Range("E31").Value = '  -3.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08708'
$ws.Range("E32").Value = '  -0.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04813'
$ws.Range("E33").Value = '  -1.62%  '

$ws.Range("E34").Value = '  -3.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.841'
$ws.Range("E35").Value = '  -0.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6849'
$ws.Range("E36").Value = '  -7.26%  '

$ws.Range("E37").Value = '  -5.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01806'
$ws.Range("E38").Value = '  -2.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.203'
$ws.Range("E39").Value = '  -8.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4914'
$ws.Range("E40").Value = '  -5.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '113.31'
$ws.Range("E41").Value = '  -2.28%  '

$ws.Range("E42").Value = '  -6.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.127'
$ws.Range("E43").Value = '  -2.28%  '

$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.748'
$ws.Range("E45").Value = '  -3.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4243'
$ws.Range("E46").Value = '  -6.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1268'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.136'
$ws.Range("E48").Value = '  -2.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05889'
$ws.Range("E49").Value = '  -0.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.93'
$ws.Range("E50").Value = '  -3.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.418'
$ws.Range("E51").Value = '  -5.77%  '
